$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "77.123.69"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").Value = "2.972.71"
$ws.Range("E3").Value = "  +3.53%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "200.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "598.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -0.26%  "

$ws.Range("E9").Value = "  +2.33%  "

$ws.Range("D10").Value = "2.962.15"
$ws.Range("E10").Value = "  +3.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.444"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +13.41%  "

$ws.Range("E12").Value = "  +0.50%  "

$ws.Range("D13").Value = "3.524.64"
$ws.Range("E13").Value = "  +3.81%  "

$ws.Range("E14").Value = "  -0.35%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "77.024.04"
$ws.Range("E15").Value = "  +1.16%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000190"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "2.952.65"
$ws.Range("E18").Value = "  +2.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.44%  "

$ws.Range("E22").Value = "  +5.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.64%  "

$ws.Range("E25").Value = "  +2.89%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("E29").Value = "  +3.41%  "

$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "500.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.401"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.42%  "

$ws.Range("B37").Value = "Cronos"
$ws.Range("C37").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.113"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +21.97%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "165.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.13%  "

$ws.Range("E41").Value = "  -5.08%  "

$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "180.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.03%  "

$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("E47").Value = "  -2.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.595"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.35%  "

